$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before the old column D (Terms Typically Offered),
# shifting it to column G.
$ws.Range("D1:F1").EntireColumn.Insert()

# Header row
$ws.Range("D1").Value2 = "Corequisites"
$ws.Range("E1").Value2 = "Concurrent"
$ws.Range("F1").Value2 = "Recommended"

# Fill in D/E/F with "NA" for all data rows, and carry G (previous D) over
# (G already holds old "Terms Typically Offered" values post-insert).
$ws.Range("D2").Value2 = "NA"
$ws.Range("E2").Value2 = "NA"
$ws.Range("F2").Value2 = "NA"

$ws.Range("D3").Value2 = "NA"
$ws.Range("E3").Value2 = "NA"
$ws.Range("F3").Value2 = "NA"

$ws.Range("D4").Value2 = "NA"
$ws.Range("E4").Value2 = "NA"
$ws.Range("F4").Value2 = "NA"

$ws.Range("D5").Value2 = "NA"
$ws.Range("E5").Value2 = "NA"
$ws.Range("F5").Value2 = "NA"

$ws.Range("D6").Value2 = "NA"
$ws.Range("E6").Value2 = "NA"
$ws.Range("F6").Value2 = "NA"

$ws.Range("D7").Value2 = "NA"
$ws.Range("E7").Value2 = "NA"
$ws.Range("F7").Value2 = "NA"

$ws.Range("D8").Value2 = "NA"
$ws.Range("E8").Value2 = "NA"
$ws.Range("F8").Value2 = "NA"

$ws.Range("D9").Value2 = "NA"
$ws.Range("E9").Value2 = "NA"
$ws.Range("F9").Value2 = "NA"

$ws.Range("D10").Value2 = "NA"
$ws.Range("E10").Value2 = "NA"
$ws.Range("F10").Value2 = "NA"

$ws.Range("D11").Value2 = "NA"
$ws.Range("E11").Value2 = "NA"
$ws.Range("F11").Value2 = "one of the BIO 327, BOT 313, BOT 326, MSCI 300, NR 305, or NR 306."
$ws.Range("G11").Value2 = "W "

$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").Value2 = "NA"
$ws.Range("F12").Value2 = "NA"

$ws.Range("D13").Value2 = "NA"
$ws.Range("E13").Value2 = "NA"
$ws.Range("F13").Value2 = "NA"

$ws.Range("D14").Value2 = "NA"
$ws.Range("E14").Value2 = "NA"
$ws.Range("F14").Value2 = "NA"

$ws.Range("D15").Value2 = "NA"
$ws.Range("E15").Value2 = "NA"
$ws.Range("F15").Value2 = "NA"

$ws.Range("D16").Value2 = "NA"
$ws.Range("E16").Value2 = "NA"
$ws.Range("F16").Value2 = "NA"

$ws.Range("D17").Value2 = "NA"
$ws.Range("E17").Value2 = "NA"
$ws.Range("F17").Value2 = "NA"

$ws.Range("D18").Value2 = "NA"
$ws.Range("E18").Value2 = "NA"
$ws.Range("F18").Value2 = "NA"

$ws.Range("D19").Value2 = "NA"
$ws.Range("E19").Value2 = "NA"
$ws.Range("F19").Value2 = "NA"

$ws.Range("D20").Value2 = "NA"
$ws.Range("E20").Value2 = "NA"
$ws.Range("F20").Value2 = "NA"

$ws.Range("D21").Value2 = "NA"
$ws.Range("E21").Value2 = "NA"
$ws.Range("F21").Value2 = "NA"

$ws.Range("D22").Value2 = "NA"
$ws.Range("E22").Value2 = "NA"
$ws.Range("F22").Value2 = "NA"

$ws.Range("D23").Value2 = "NA"
$ws.Range("E23").Value2 = "NA"
$ws.Range("F23").Value2 = "NA"

$ws.Range("D24").Value2 = "NA"
$ws.Range("E24").Value2 = "NA"
$ws.Range("F24").Value2 = "NA"

$ws.Range("D25").Value2 = "NA"
$ws.Range("E25").Value2 = "NA"
$ws.Range("F25").Value2 = "NA"

$ws.Range("D26").Value2 = "NA"
$ws.Range("E26").Value2 = "NA"
$ws.Range("F26").Value2 = "NA"

$ws.Range("D27").Value2 = "NA"
$ws.Range("E27").Value2 = "NA"
$ws.Range("F27").Value2 = "NA"

# Update prerequisite text in C7 and C11 to reflect new requirement separation
$ws.Range("C7").Value2 = "one of the AEPS 120, BOT 121, or SS 120; and CHEM 124 or CHEM 127."
$ws.Range("C11").Value2 = "BOT 121 or BIO 162; CHEM 127; and SS 120 or SS 130."
